$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 954.15
$ws.Range("I17").Value = 862.25
$ws.Range("K17").Value = 2586.75
$ws.Range("M17").Value = -2418.75
$ws.Range("H98").Value = 3739.2
$ws.Range("I98").Value = 4199
$ws.Range("J98").Value = 1900
$ws.Range("K98").Value = 4199
$ws.Range("L98").Value = 1900
$ws.Range("M98").Value = -2701
$ws.Range("N98").Value = -4896
$ws.Range("H112").Value = 1816.1666
$ws.Range("J112").Value = 1742.5714
$ws.Range("L112").Value = 5227.7142
$ws.Range("N112").Value = -7443.7142
$ws.Range("H113").Value = 3578.913
$ws.Range("I113").Value = 3814
$ws.Range("J113").Value = 3322.4546
$ws.Range("K113").Value = 3814
$ws.Range("L113").Value = 3322.4546
$ws.Range("M113").Value = -560
$ws.Range("N113").Value = -9830.454600000001
$ws.Range("H116").Value = 20587.215
$ws.Range("I116").Value = 23090.625
$ws.Range("J116").Value = 17249.334
$ws.Range("K116").Value = 23090.625
$ws.Range("L116").Value = 17249.334
$ws.Range("M116").Value = -19648.625
$ws.Range("N116").Value = -24133.334
$ws.Range("H121").Value = 4850
$ws.Range("J121").Value = 4850
$ws.Range("L121").Value = 14550
$ws.Range("N121").Value = -18044
$ws.Range("H122").Value = 3739.2
$ws.Range("I122").Value = 4199
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 12597
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -10147
$ws.Range("N122").Value = -10600
$ws.Range("H129").Value = 2116.5
$ws.Range("J129").Value = 2175
$ws.Range("L129").Value = 6525
$ws.Range("N129").Value = -16525

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 994828.9
$ws.Range("I32").Value = 1049007.9
$ws.Range("K32").Value = 1049007.9
$ws.Range("M32").Value = -1048720.9
$ws.Range("H61").Value = 2278273.8
$ws.Range("I61").Value = 6391.6553
$ws.Range("K61").Value = 6391.6553
$ws.Range("M61").Value = -6179.6553
$ws.Range("H74").Value = 2541458
$ws.Range("I74").Value = 2790404
$ws.Range("K74").Value = 2790404
$ws.Range("M74").Value = -2789530
$ws.Range("H77").Value = 2541458
$ws.Range("I77").Value = 2790404
$ws.Range("K77").Value = 13952020
$ws.Range("M77").Value = -13947652
$ws.Range("H97").Value = 43479220
$ws.Range("I97").Value = 980.8421
$ws.Range("K97").Value = 980.8421
$ws.Range("M97").Value = -484.8421
$ws.Range("H136").Value = 2278273.8
$ws.Range("I136").Value = 6391.6553
$ws.Range("K136").Value = 19174.9659
$ws.Range("M136").Value = -16624.9659

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 68440.42
$ws.Range("J20").Value = 1518.8
$ws.Range("L20").Value = 1518.8
$ws.Range("N20").Value = -2012.8
$ws.Range("H86").Value = 2690.8965
$ws.Range("I86").Value = 2277.818
$ws.Range("K86").Value = 2277.818
$ws.Range("M86").Value = -1154.818
$ws.Range("H89").Value = 2690.8965
$ws.Range("I89").Value = 2277.818
$ws.Range("K89").Value = 11389.09
$ws.Range("M89").Value = -5773.09
$ws.Range("H105").Value = 6137.1177
$ws.Range("I105").Value = 2986.0833
$ws.Range("J105").Value = 13699.6
$ws.Range("K105").Value = 2986.0833
$ws.Range("L105").Value = 13699.6
$ws.Range("M105").Value = -1239.0833
$ws.Range("N105").Value = -17193.6

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 36665
$ws.Range("I59").Value = 29997.5
$ws.Range("J59").Value = 50000
$ws.Range("K59").Value = 29997.5
$ws.Range("L59").Value = 50000
$ws.Range("M59").Value = -28852.5
$ws.Range("N59").Value = -52290
$ws.Range("H86").Value = 33150.617
$ws.Range("I86").Value = 66442
$ws.Range("K86").Value = 66442
$ws.Range("M86").Value = -65319
$ws.Range("H89").Value = 33150.617
$ws.Range("I89").Value = 66442
$ws.Range("K89").Value = 332210
$ws.Range("M89").Value = -326594
$ws.Range("H105").Value = 30598.428
$ws.Range("I105").Value = 17364.834
$ws.Range("K105").Value = 17364.834
$ws.Range("M105").Value = -15617.834
$ws.Range("H107").Value = 789.3333
$ws.Range("I107").Value = 616.4167
$ws.Range("J107").Value = 1135.1666
$ws.Range("K107").Value = 616.4167
$ws.Range("L107").Value = 1135.1666
$ws.Range("M107").Value = 1303.5833
$ws.Range("N107").Value = -4975.1666
$ws.Range("H115").Value = 44996.668
$ws.Range("J115").Value = 44996.668
$ws.Range("L115").Value = 44996.668
$ws.Range("N115").Value = -47346.668

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 4548.8
$ws.Range("I94").Value = 4548.8
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 13646.4
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -12970.4
$ws.Range("H107").Value = 927.75
$ws.Range("I107").Value = 1317.8
$ws.Range("J107").Value = 649.1429000000001
$ws.Range("K107").Value = 3953.4
$ws.Range("L107").Value = 1947.4287
$ws.Range("M107").Value = -2033.4
$ws.Range("N107").Value = -5787.4287
$ws.Range("H119").Value = 4252.5713
$ws.Range("I119").Value = 2461.3333
$ws.Range("J119").Value = 15000
$ws.Range("K119").Value = 7383.999899999999
$ws.Range("L119").Value = 45000
$ws.Range("M119").Value = -2545.999899999999
$ws.Range("N119").Value = -54676
$ws.Range("H124").Value = 8188
$ws.Range("I124").Value = 8188
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 24564
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -19654
$ws.Range("N94").ClearContents()
$ws.Range("N124").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15968
$ws.Range("I70").Value = 20545.334
$ws.Range("J70").Value = 6813.3335
$ws.Range("K70").Value = 20545.334
$ws.Range("L70").Value = 6813.3335
$ws.Range("M70").Value = -20275.334
$ws.Range("N70").Value = -7353.3335
$ws.Range("H73").Value = 15968
$ws.Range("I73").Value = 20545.334
$ws.Range("J73").Value = 6813.3335
$ws.Range("K73").Value = 20545.334
$ws.Range("L73").Value = 6813.3335
$ws.Range("M73").Value = -19609.334
$ws.Range("N73").Value = -8685.333500000001
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("M83").Value = -5008
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3256.4546
$ws.Range("I22").Value = 1986
$ws.Range("J22").Value = 3982.4285
$ws.Range("K22").Value = 1986
$ws.Range("L22").Value = 3982.4285
$ws.Range("M22").Value = -1691
$ws.Range("N22").Value = -4572.4285
$ws.Range("H27").Value = 3256.4546
$ws.Range("I27").Value = 1986
$ws.Range("J27").Value = 3982.4285
$ws.Range("K27").Value = 1986
$ws.Range("L27").Value = 3982.4285
$ws.Range("M27").Value = -1879
$ws.Range("N27").Value = -4196.4285
$ws.Range("H122").Value = 4389.5
$ws.Range("J122").Value = 4294.3335
$ws.Range("L122").Value = 12883.0005
$ws.Range("N122").Value = -17783.0005
$ws.Range("H136").Value = 9262385
$ws.Range("I136").Value = 6253118
$ws.Range("K136").Value = 18759354
$ws.Range("M136").Value = -18756804

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 113630.11
$ws.Range("I4").Value = 144809.28
$ws.Range("K4").Value = 144809.28
$ws.Range("M4").Value = -144696.28
$ws.Range("H126").Value = 1153.3529
$ws.Range("I126").Value = 800.5
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 2401.5
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = 68.5
$ws.Range("N126").Value = -13340
$ws.Range("H136").Value = 4867169
$ws.Range("I136").Value = 2293525
$ws.Range("J136").Value = 21166916
$ws.Range("K136").Value = 6880575
$ws.Range("L136").Value = 63500748
$ws.Range("M136").Value = -6878025
$ws.Range("N136").Value = -63505848
